$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update "Förändrad" (last-changed) date column C for rows 2-5
# from 45184 (2023-09-15) to 45185 (2023-09-16)
$ws.Range("C2:C5").Value = 45185
